# Updates the "cryptos" price-tracker sheet with a fresh snapshot of
# prices / 1h-volume figures (and, for rows 39-40, a re-ranking swap
# between Aptos and InternetComputer(DFINITY)).
#
# The sheet stores every data cell as literal text (prices such as
# "1.002" or "28.315.93" are NOT numbers - some even contain two dots -
# and the Volume(1h) column is padded, e.g. "  -0.46%  "). Writing a
# numeric-looking string straight into `.Value` makes Excel silently
# reinterpret it as a real number (and reformat it, e.g. "1.003" ->
# 1.0029999999999999), which would corrupt the data. Prefixing the
# string with a leading apostrophe forces Excel to store it verbatim as
# text; re-applying the built-in "Normal" style afterwards clears the
# quote-prefix flag that the apostrophe trick leaves behind so the
# cell's style stays identical to how it started (no stray "treat as
# text" marker left on the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> {Col = NewValue} update table, columns: B=Coin, C=Link, D=Price, E=Volume(1h)
$updates = @(
    @{ Row=2; D='28.319.82'; E='  -0.52%  ' }
    @{ Row=3; D='1.807.49'; E='  -0.94%  ' }
    @{ Row=4; D='1.003'; E='  -0.17%  ' }
    @{ Row=5; D='313.39'; E='  -1.16%  ' }
    @{ Row=7; D='0.5150'; E='  -0.25%  ' }
    @{ Row=8; D='0.3970'; E='  +2.83%  ' }
    @{ Row=9; D='0.07841'; E='  -5.45%  ' }
    @{ Row=10; D='1.111'; E='  -1.07%  ' }
    @{ Row=11; D='41.02'; E='  -2.08%  ' }
    @{ Row=12; D='6.321'; E='  -0.94%  ' }
    @{ Row=13; E='  -0.20%  ' }
    @{ Row=14; D='20.42'; E='  -3.17%  ' }
    @{ Row=15; D='1.812.14'; E='  -0.67%  ' }
    @{ Row=16; D='7.309' }
    @{ Row=17; D='92.62'; E='  -1.64%  ' }
    @{ Row=18; D='0.00001084'; E='  -3.37%  ' }
    @{ Row=19; D='0.06562'; E='  -1.13%  ' }
    @{ Row=20; E='  -0.12%  ' }
    @{ Row=21; D='17.28'; E='  -2.86%  ' }
    @{ Row=22; D='6.013'; E='  -0.67%  ' }
    @{ Row=23; D='28.363.74'; E='  -0.44%  ' }
    @{ Row=24; D='11.12'; E='  -3.34%  ' }
    @{ Row=25; D='2.228'; E='  -0.82%  ' }
    @{ Row=26; D='161.09'; E='  +1.07%  ' }
    @{ Row=27; D='20.51'; E='  -2.79%  ' }
    @{ Row=28; D='2.020.71'; E='  -0.66%  ' }
    @{ Row=29; D='2.415'; E='  +0.37%  ' }
    @{ Row=30; D='127.83'; E='  +1.42%  ' }
    @{ Row=31; D='0.1103'; E='  -0.59%  ' }
    @{ Row=32; D='1.062'; E='  -2.79%  ' }
    @{ Row=33; D='3.666'; E='  -0.50%  ' }
    @{ Row=34; D='5.572'; E='  -2.77%  ' }
    @{ Row=35; D='0.07190'; E='  -4.65%  ' }
    @{ Row=36; D='9.128'; E='  +4.11%  ' }
    @{ Row=37; D='0.02358'; E='  -0.14%  ' }
    @{ Row=38; D='0.2190'; E='  -1.62%  ' }
    @{ Row=39; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.053'; E='  -3.67%  ' }
    @{ Row=40; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='11.58'; E='  -4.69%  ' }
    @{ Row=41; D='0.6193'; E='  -3.08%  ' }
    @{ Row=42; D='1.001'; E='  -0.18%  ' }
    @{ Row=43; E='  -2.34%  ' }
    @{ Row=44; D='13.19'; E='  -2.70%  ' }
    @{ Row=45; D='0.5988'; E='  -3.28%  ' }
    @{ Row=46; D='3.747'; E='  -1.40%  ' }
    @{ Row=47; E='  -6.73%  ' }
    @{ Row=48; D='125.46'; E='  -1.89%  ' }
    @{ Row=49; D='1.217'; E='  +0.89%  ' }
    @{ Row=50; D='1.923'; E='  -4.14%  ' }
    @{ Row=51; D='0.06830'; E='  -1.88%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in 'B', 'C', 'D', 'E') {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Cells.Item($row, [int][char]$col - [int][char]'A' + 1)
            $cell.Value = "'" + $u[$col]
            $cell.Style = "Normal"
        }
    }
}
